$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rushing")
$ws2 = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------------
# Sheet "Rushing": update stats produced by simulating the Wild Card round
# ---------------------------------------------------------------------------

# Row 3 - T.Huntley
$ws1.Cells.Item(3, 3).Value = 15
$ws1.Cells.Item(3, 4).Value = 16
$ws1.Cells.Item(3, 5).Value = 13
$ws1.Cells.Item(3, 6).Value = 7

# Row 5 - T.Williams
$ws1.Cells.Item(5, 3).Value = 22
$ws1.Cells.Item(5, 6).Value = 4

# Row 6 - L.Murray
$ws1.Cells.Item(6, 3).Value = 68
$ws1.Cells.Item(6, 4).Value = 38
$ws1.Cells.Item(6, 5).Value = 12

# Row 7 - D.Freeman
$ws1.Cells.Item(7, 3).Value = 78

# Row 11 - D.Duvernay
$ws1.Cells.Item(11, 3).Value = 3

# New row 12 - M.Andrews joins the rushing log
$ws1.Cells.Item(11, 1).Copy()
$ws1.Cells.Item(12, 1).PasteSpecial(-4122)
$ws1.Cells.Item(12, 1).Value = 10
$ws1.Cells.Item(12, 2).Value = "M.Andrews"
$ws1.Cells.Item(12, 3).Value = 0
$ws1.Cells.Item(12, 4).Value = 0
$ws1.Cells.Item(12, 5).Value = 1
$ws1.Cells.Item(12, 6).Value = 1

# ---------------------------------------------------------------------------
# Sheet "Receiving": update stats produced by simulating the Wild Card round
# ---------------------------------------------------------------------------

# Row 4 - D.Freeman
$ws2.Cells.Item(4, 3).Value = 41
$ws2.Cells.Item(4, 4).Value = 34
$ws2.Cells.Item(4, 7).Value = 9
$ws2.Cells.Item(4, 8).Value = 7

# Row 6 - M.Brown
$ws2.Cells.Item(6, 3).Value = 101
$ws2.Cells.Item(6, 4).Value = 76
$ws2.Cells.Item(6, 5).Value = 44
$ws2.Cells.Item(6, 7).Value = 16

# Row 7 - S.Watkins
$ws2.Cells.Item(7, 3).Value = 36

# Row 10 - T.Wallace
$ws2.Cells.Item(10, 3).Value = 5

# Row 11 - R.Bateman
$ws2.Cells.Item(11, 3).Value = 51
$ws2.Cells.Item(11, 4).Value = 37

# Row 13 - M.Andrews
$ws2.Cells.Item(13, 3).Value = 118
$ws2.Cells.Item(13, 4).Value = 87
$ws2.Cells.Item(13, 5).Value = 36
$ws2.Cells.Item(13, 6).Value = 20
$ws2.Cells.Item(13, 7).Value = 22

# Row 15 - E.Tomlinson
$ws2.Cells.Item(15, 3).Value = 2
$ws2.Cells.Item(15, 7).Value = 1

# Row 16 - J.Oliver
$ws2.Cells.Item(16, 3).Value = 13
$ws2.Cells.Item(16, 4).Value = 9
